# Restored from revision #7cb052631180c5404a8538e530db3f9ca40d266c.TEST
# Author: admin. Type: SAVE.
#
# Functional change: cell C10 on the active ("Rules") sheet changes
# from 18 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
